# 27 Mayis (May 27, 2020) verileri eklendi / add 27 May 2020 data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds an Excel Table ("Table3") over A1:E76 that must grow by
# one row to cover the newly appended data (A1:E77).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E77"))

# Append the new day's figures (date serial 43978 = 2020-05-27) in row 77:
# date, test, case, death, recovered.
$ws.Cells.Item(77, 1).Value = 43978
$ws.Cells.Item(77, 2).Value = 21043
$ws.Cells.Item(77, 3).Value = 1035
$ws.Cells.Item(77, 4).Value = 34
$ws.Cells.Item(77, 5).Value = 1286

# Match Excel's natural post-entry selection (one row below the new data).
$ws.Range("B78").Select()
